$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.271.09"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.361.87"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.09"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.02"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.361.92"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.629"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.92"
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.13"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.882.83"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.09"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.354.17"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "65.213.57"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.72"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.989"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "474.64"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.91"
$ws.Range("E23").Value = "  -6.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.92"
$ws.Range("E24").Value = "  +3.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.06"
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.02"
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.87"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.51"
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.62"
$ws.Range("E29").Value = "  -3.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.01"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.48"
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.34"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "571.45"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.107"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.62"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.35"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.370"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0732"
$ws.Range("E41").Value = "  -3.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.083.59"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.77"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0412"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.133"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.17"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.996"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.23"
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.56"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.40"
$ws.Range("E51").Value = "  +0.07%  "
